# Auto-update Cloudflare export:
# For every DNS record row whose "type" (column C) is "CNAME",
# set its "settings" (column H) from the empty dict "{}" to
# "{'flatten_cname': False}".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$updated = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $recordType = $ws.Cells.Item($r, 3).Value2
    if ($recordType -eq "CNAME") {
        $settingsCell = $ws.Cells.Item($r, 8)
        if ($settingsCell.Value2 -eq "{}") {
            $settingsCell.Value = "{'flatten_cname': False}"
            $updated = $updated + 1
        }
    }
}

Write-Output "Updated $updated CNAME settings cells"
